$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 13)
    $cell.Value = $cell.Value2 * 100
}

$ws.Range("M2:M51").NumberFormat = "0.00"
